# Daily rollover update: advance the "days remaining" counters by one day.
#
# Column layout (Sheet1):
#   A = row number, B = shop name, C = address,
#   D = total days (总天), E = days remaining (剩余), F = start date (开始时间, yyyymmdd),
#   G/H/I = notes.
#
# For every data row:
#   - if the remaining-days counter (E) is already at 1 (i.e. it would hit 0
#     today), the cycle restarts: E is reset to the total-days value (D) and
#     the start date (F) is rolled to the new "today" (2026-01-06).
#   - otherwise E simply decrements by 1 and F is left untouched.
#
# Rows whose start-date cell (F) isn't a well-formed 8-digit yyyymmdd value
# (e.g. a corrupted "202510929") are left alone entirely, same as upstream.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newToday = 20260106
$lastRow = $ws.UsedRange.Rows.Count
if ($lastRow -lt 2) {
    $lastRow = 99
}

for ($r = 2; $r -le $lastRow; $r++) {
    $totalDays = $ws.Cells.Item($r, 4).Value2
    $remaining = $ws.Cells.Item($r, 5).Value2
    $startDate = $ws.Cells.Item($r, 6).Value2

    if ($remaining -eq $null -or $totalDays -eq $null) {
        continue
    }

    $dateText = [string]$startDate
    if ($dateText.Length -ne 8) {
        # Malformed start date (not yyyymmdd) - skip, matching the source data.
        continue
    }

    if ($remaining -eq 1) {
        $ws.Cells.Item($r, 5).Value2 = $totalDays
        $ws.Cells.Item($r, 6).Value2 = $newToday
    } else {
        $ws.Cells.Item($r, 5).Value2 = $remaining - 1
    }
}
